$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the arbitrary stimuli labels (column B-I of row 2) with the
# standardised set used by Leader and Barnes-Holmes's (2001b)
$ws.Range("B2").Value = "ZID"
$ws.Range("C2").Value = "VEK"
$ws.Range("D2").Value = "YIM"
$ws.Range("E2").Value = "PAF"
$ws.Range("F2").Value = "ROG"
$ws.Range("G2").Value = "MAU"
$ws.Range("H2").Value = "JOM"
$ws.Range("I2").Value = "DAX"

# Set an (unused) cell's font further down the sheet, which is how the
# extra LiberationSerif font / style entries ended up in the workbook
$ws.Range("B11").Font.Name = "LiberationSerif"

# Leave the selection where the new values were entered
[void]$ws.Range("A2:I2").Select()
